$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.2648423062635428
$ws.Cells.Item(2, 9).Value = 0.3999845054128732
$ws.Cells.Item(2, 10).Value = 0.1667020416781071
$ws.Cells.Item(2, 11).Value = 0.9893093490273974
$ws.Cells.Item(2, 15).Value = 0.2314577481323485
$ws.Cells.Item(3, 7).Value = 0.7367725734787592
$ws.Cells.Item(3, 9).Value = 0.8067698848656968
$ws.Cells.Item(3, 10).Value = 0.7013610185265239
$ws.Cells.Item(3, 11).Value = 1.079606318662967
$ws.Cells.Item(3, 13).Value = 1.535090727082321
$ws.Cells.Item(3, 15).Value = 0.5874117869445946
$ws.Cells.Item(4, 7).Value = 0.3172829080795564
$ws.Cells.Item(4, 9).Value = 0.4556921325965645
$ws.Cells.Item(4, 10).Value = 0.2329149865727795
$ws.Cells.Item(4, 11).Value = 1.006599966677753
$ws.Cells.Item(4, 13).Value = 0.8296464212590471
$ws.Cells.Item(4, 15).Value = 0.3011757366550303
$ws.Cells.Item(5, 7).Value = 1.047167857446734
$ws.Cells.Item(5, 9).Value = 0.8542726440854786
$ws.Cells.Item(5, 10).Value = 1.486231553615516
$ws.Cells.Item(5, 11).Value = 0.8598011502695897
$ws.Cells.Item(5, 13).Value = 0.3991767594233399
$ws.Cells.Item(5, 15).Value = 1.349922214039842
$ws.Cells.Item(6, 7).Value = 0.990122340778822
$ws.Cells.Item(6, 9).Value = 0.7849234254131709
$ws.Cells.Item(6, 10).Value = 0.491954426243003
$ws.Cells.Item(6, 11).Value = 0.8486782587928934
$ws.Cells.Item(6, 13).Value = 0.3703395202782027
$ws.Cells.Item(6, 15).Value = 0.2988408563168703
$ws.Cells.Item(7, 7).Value = 0.4488583611760173
$ws.Cells.Item(7, 9).Value = 0.3576788438215815
$ws.Cells.Item(7, 10).Value = 0.2735505756474028
$ws.Cells.Item(7, 11).Value = 0.8957167478771801
$ws.Cells.Item(7, 13).Value = 0.2605979402896592
$ws.Cells.Item(7, 15).Value = 0.3672874443582904
$ws.Cells.Item(8, 7).Value = 0.9062454332425819
$ws.Cells.Item(8, 9).Value = 0.7362623963781584
$ws.Cells.Item(8, 10).Value = 0.3023720181096924
$ws.Cells.Item(8, 11).Value = 0.8418244561768575
$ws.Cells.Item(8, 13).Value = 0.3392403909474527
$ws.Cells.Item(8, 15).Value = 0.3484605790162927
$ws.Cells.Item(9, 7).Value = 1.237113586789377
$ws.Cells.Item(9, 9).Value = 0.9549235379643154
$ws.Cells.Item(9, 10).Value = 0.363242646818582
$ws.Cells.Item(9, 11).Value = 0.8380172508133381
$ws.Cells.Item(9, 13).Value = 0.5358669738276819
$ws.Cells.Item(9, 15).Value = 0.4564606114509621
$ws.Cells.Item(10, 7).Value = 0.2272709963872031
$ws.Cells.Item(10, 9).Value = 0.3505290519841819
$ws.Cells.Item(10, 10).Value = 1.105161635716102
$ws.Cells.Item(10, 11).Value = 0.9744048986904357
$ws.Cells.Item(10, 13).Value = 0.5866892550215309
$ws.Cells.Item(10, 15).Value = 0.8871267945290771
$ws.Cells.Item(11, 7).Value = 1.740248933809657
$ws.Cells.Item(11, 9).Value = 1.472665161823985
$ws.Cells.Item(11, 10).Value = 1.470261269805833
$ws.Cells.Item(11, 11).Value = 0.9237447867456334
$ws.Cells.Item(11, 13).Value = 1.006566124253472
$ws.Cells.Item(11, 15).Value = 1.636274910450251
$ws.Cells.Item(12, 7).Value = 2.808236504563006
$ws.Cells.Item(12, 9).Value = 2.43911637503443
$ws.Cells.Item(12, 10).Value = 1.469662336551587
$ws.Cells.Item(12, 11).Value = 1.707287010441392
$ws.Cells.Item(12, 13).Value = 2.066582927066592
$ws.Cells.Item(12, 15).Value = 1.257223957844785
$ws.Cells.Item(13, 7).Value = 2.96375479385902
$ws.Cells.Item(13, 9).Value = 2.608351823819266
$ws.Cells.Item(13, 10).Value = 1.036435725319554
$ws.Cells.Item(13, 11).Value = 1.780065299615569
$ws.Cells.Item(13, 13).Value = 2.223288467225127
$ws.Cells.Item(13, 15).Value = 0.7801894478578744
$ws.Cells.Item(14, 7).Value = 0.5293721803614317
$ws.Cells.Item(14, 10).Value = 0.3780563917137483
$ws.Cells.Item(14, 11).Value = 0.6705288111468155
$ws.Cells.Item(15, 7).Value = 0.3214587559052511
$ws.Cells.Item(15, 10).Value = 0.2620851585266541
$ws.Cells.Item(15, 11).Value = 0.7011247130273096
$ws.Cells.Item(15, 13).Value = 0.790585644941196
$ws.Cells.Item(16, 7).Value = 0.6195532464453075
$ws.Cells.Item(16, 10).Value = 0.5750369945155545
$ws.Cells.Item(16, 11).Value = 0.8723360262831423
$ws.Cells.Item(16, 13).Value = 0.2882802641105014
$ws.Cells.Item(17, 7).Value = 0.5044473022233332
$ws.Cells.Item(17, 10).Value = 0.2202018135256127
$ws.Cells.Item(17, 11).Value = 0.8551202753362673
$ws.Cells.Item(17, 13).Value = 0.2677362010018021
$ws.Cells.Item(18, 7).Value = 0.1769659046176615
$ws.Cells.Item(18, 9).Value = 0.2847458844189225
$ws.Cells.Item(18, 10).Value = 0.231952375889422
$ws.Cells.Item(18, 11).Value = 0.75329983683942
$ws.Cells.Item(18, 13).Value = 0.456913988622148
$ws.Cells.Item(19, 7).Value = 0.4603100215000683
$ws.Cells.Item(19, 9).Value = 0.5425801482383638
$ws.Cells.Item(19, 10).Value = 0.4295909304226662
$ws.Cells.Item(19, 11).Value = 0.6833154897163575
$ws.Cells.Item(19, 13).Value = 1.00899767024709
$ws.Cells.Item(20, 7).Value = 0.1782568086682427
$ws.Cells.Item(20, 9).Value = 0.2827058914318142
$ws.Cells.Item(20, 10).Value = 0.2703431646217884
$ws.Cells.Item(20, 11).Value = 0.7512184662674252
$ws.Cells.Item(20, 13).Value = 0.4701996238443251
$ws.Cells.Item(21, 7).Value = 0.1957384692117254
$ws.Cells.Item(21, 9).Value = 0.2931407093753178
$ws.Cells.Item(21, 10).Value = 0.2233767533978894
$ws.Cells.Item(21, 11).Value = 0.7861187593041835
$ws.Cells.Item(21, 13).Value = 0.3320203296529319
$ws.Cells.Item(22, 7).Value = 0.2995617203586911
$ws.Cells.Item(22, 9).Value = 0.3183782736431794
$ws.Cells.Item(22, 10).Value = 0.3338913416917921
$ws.Cells.Item(22, 11).Value = 0.8194224785569194
$ws.Cells.Item(22, 13).Value = 0.2670639643961904
$ws.Cells.Item(23, 7).Value = 0.3018427706456829
$ws.Cells.Item(23, 9).Value = 0.3019466076076006
$ws.Cells.Item(23, 10).Value = 0.2303339987425469
$ws.Cells.Item(23, 11).Value = 0.8203472329602192
$ws.Cells.Item(23, 13).Value = 0.2666843675498428
$ws.Cells.Item(24, 7).Value = 0.7859439470026772
$ws.Cells.Item(24, 9).Value = 0.8647302407095083
$ws.Cells.Item(24, 10).Value = 0.8053580720500954
$ws.Cells.Item(24, 11).Value = 0.6725266679115
$ws.Cells.Item(24, 13).Value = 1.406119705266911
$ws.Cells.Item(25, 7).Value = 0.1676900696445014
$ws.Cells.Item(25, 9).Value = 0.2781656602250006
$ws.Cells.Item(25, 10).Value = 0.4387373481867253
$ws.Cells.Item(25, 11).Value = 0.7693459925171294
$ws.Cells.Item(25, 13).Value = 0.3896174153040006
